# Auto-generated Excel COM-interop script
# Applies the "Linea 141" schedule-refresh edit (commit: "Horarios actualizados Linea 141 - 368")
# to the three worksheets of the workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "LP1912" — full table grows from 29 data rows to 41 data rows.
#   * Header timestamps (A2/A3) refreshed.
#   * 12 new rows are inserted before the old row 27 so the table grows from
#     rows 27-34 (8 rows) to rows 27-46 (20 rows).
#   * All 20 rows (27-46) are then (re)written with the new/after data so the
#     existing rows that shifted down also land on their final values.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2,1).Value = "Última actualización: 05:49:10"
$ws1.Cells.Item(3,1).Value = "Total filas: 41"

$ws1.Range("A27:A38").EntireRow.Insert()

$data1 = @(
    @("05:49:10", "06:30", "86_EST CHICA-ESC AGRARIA", 41, "LP1912"),
    @("04:40:48", "06:31", "16_SANTA ANA", 111, "LP1912"),
    @("04:54:17", "06:44", "225_C ROCA-H SUR", 110, "LP1912"),
    @("04:54:17", "06:46", "215C_EL PATO", 112, "LP1912"),
    @("05:49:10", "06:47", "215C_EL PATO", 58, "LP1912"),
    @("05:19:24", "06:59", "14_ABASTO", 100, "LP1912"),
    @("05:49:10", "07:00", "14_ABASTO", 71, "LP1912"),
    @("05:49:10", "07:05", "23_HERNANDEZ", 76, "LP1912"),
    @("05:19:24", "07:05", "15_ABASTO", 106, "LP1912"),
    @("05:19:24", "07:07", "225_GOMEZ", 108, "LP1912"),
    @("05:19:24", "07:11", "215A_EL PATO", 112, "LP1912"),
    @("05:19:24", "07:15", "11_ETCHEVERRY", 116, "LP1912"),
    @("05:49:10", "07:16", "11_ETCHEVERRY", 87, "LP1912"),
    @("05:49:10", "07:21", "26_HERNANDEZ", 92, "LP1912"),
    @("05:49:10", "07:32", "11_ETCHEVERRY", 103, "LP1912"),
    @("05:49:10", "07:32", "16_SANTA ANA", 103, "LP1912"),
    @("05:49:10", "07:32", "84_COLONIA URQUIZA-ESC 49", 103, "LP1912"),
    @("05:49:10", "07:37", "27_EL RETIRO", 108, "LP1912"),
    @("05:49:10", "07:39", "10_OLMOS", 110, "LP1912"),
    @("05:49:10", "07:48", "14_ABASTO", 119, "LP1912")
)

$r = 27
foreach ($row in $data1) {
    for ($c = 1; $c -le 5; $c++) {
        $ws1.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet 2: "LP1912-215" — grows from 7 to 8 data rows.
#   * Header timestamps (A2/A3) refreshed.
#   * One row is inserted at row 12 (pushing the former row 12 down to 13)
#     and populated with the new data.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2,1).Value = "Última actualización: 05:49:10"
$ws2.Cells.Item(3,1).Value = "Total filas: 8"

$ws2.Rows.Item(12).Insert()

$row12 = @("05:49:10", "06:47", "215C_EL PATO", 58, "LP1912")

for ($c = 1; $c -le 5; $c++) {
    $ws2.Cells.Item(12, $c).Value = $row12[$c - 1]
}

# ---------------------------------------------------------------------------
# Sheet 3: "6203-6173" — grows from 7 to 8 data rows.
#   * Header timestamps (A2/A3) refreshed.
#   * One new row is appended as row 13 (no shift needed; nothing below it).
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2,1).Value = "Última actualización: 05:49:10"
$ws3.Cells.Item(3,1).Value = "Total filas: 8"

$row13 = @("05:49:10", "07:35", "215A_LA PLATA", 106, "L6173")

for ($c = 1; $c -le 5; $c++) {
    $ws3.Cells.Item(13, $c).Value = $row13[$c - 1]
}
